# Auto-generated PowerShell-style Excel COM-interop script
# Applies the cryptos.xlsx data-refresh diff (prices & % changes updated,
# plus TrustWalletToken/Cronos row order swap) cell by cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.286.91"
$ws.Range("E2").Value = "  -1.17%  "
$ws.Range("D3").Value = "2.280.78"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "265.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.52%  "
$ws.Range("E7").Value = "  -0.90%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.608"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0930"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("D15").Value = "2.624.03"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.861"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "2.278.60"
$ws.Range("E17").Value = "  -1.80%  "
$ws.Range("D18").Value = "43.215.35"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("E19").Value = "  -2.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = "  -2.58%  "
$ws.Range("E23").Value = "  -1.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "40.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.13%  "
$ws.Range("E29").Value = "  -2.32%  "
$ws.Range("E30").Value = "  -1.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0906"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.41%  "
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("E36").Value = "  -2.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.96"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("E39").Value = "  -5.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "79.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.238"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.35%  "
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("E46").Value = "  -1.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.81"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.34%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.34%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0995"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.442"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.88%  "
